# Excel and BOM updated
# Adds a new BOM line item (ON Schottky Diode / MBRM110LT1G / D2) to the
# CabinTempSenseHardware BOM sheet, right before the "Per Board" / "Order
# Total" summary rows, and extends the summary formulas to include it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old row 16 ("Per Board"/"Order Total" label
# row), which pushes that row down to 17 and the totals row down to 18.
$ws.Rows("16:16").Insert()

# Populate the new row 16 with the new part.
$ws.Range("A16").Value = "ON Schottky Diode"
$ws.Range("B16").Value = "MBRM110LT1G"
$ws.Range("C16").Value = "D2"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = 0.554
$ws.Range("G16").Value = "https://www.digikey.ca/product-detail/en/on-semiconductor/MBRM110LT1G/MBRM110LT1GOSCT-ND/917992"
$ws.Range("H16").Formula = "=E16*F16"

# Turn G16 into a real hyperlink, and make sure it keeps the same
# "Hyperlink" look used by the other Digikey links in column G.
$ws.Hyperlinks.Add($ws.Range("G16"), "https://www.digikey.ca/product-detail/en/on-semiconductor/MBRM110LT1G/MBRM110LT1GOSCT-ND/917992")
$ws.Range("G16").Style = "Hyperlink"

# Extend the running-total formulas (now on row 18) to include the new row.
$ws.Range("F18").Formula = "=D2*F2+D3*F3+D4*F4+D5*F5+D6*F6+D7*F7+D8*F8+D9*F9+D10*F10+D11*F11+D12*F12+D13*F13+D14*F14+D15*F15+D16*F16"
$ws.Range("H18").Formula = "=H2+H3+H4+H5+H6+H7+H8+H9+H10+H11+H12+H13+H14+H15+H16"

# Re-enter the H4:H10 "quantity*price" formulas as one shared formula group
# (matches the canonical save produced by Excel for this column).
$ws.Range("H4:H10").Formula = "=E4*F4"

# Restore the previously-selected cell (shifted from G26 to G22 now that
# the sheet only grew by one data row before the summary rows).
$ws.Range("G22").Select()
